# Apply crypto list updates (price/volume refresh + two coin-row swaps)
# Sheet1, per the diff: Litecoin<->ShibaInu (rows 20/21),
# InjectiveProtocol<->Aave (rows 46/47), Maker<->FraxShare (rows 48/49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "36.313.52"
Set-TextCell "E2" "  -3.12%  "
Set-TextCell "D3" "1.978.49"
Set-TextCell "E3" "  -3.82%  "
Set-TextCell "E4" "  +0.03%  "
Set-TextCell "D5" "244.39"
Set-TextCell "E5" "  -3.55%  "
Set-TextCell "E6" "  -4.50%  "
Set-TextCell "D7" "58.53"
Set-TextCell "E7" "  -12.93%  "
Set-TextCell "E8" "  +0.02%  "
Set-TextCell "D9" "0.372"
Set-TextCell "E9" "  -6.18%  "
Set-TextCell "D10" "56.96"
Set-TextCell "E10" "  -4.94%  "
Set-TextCell "D11" "0.0836"
Set-TextCell "E11" "  +8.52%  "
Set-TextCell "E12" "  -0.75%  "
Set-TextCell "D13" "23.03"
Set-TextCell "E13" "  -2.93%  "
Set-TextCell "D14" "0.856"
Set-TextCell "E14" "  -8.48%  "
Set-TextCell "D15" "2.268.55"
Set-TextCell "E15" "  -3.86%  "
Set-TextCell "D16" "13.89"
Set-TextCell "E16" "  -7.29%  "
Set-TextCell "D17" "5.42"
Set-TextCell "E17" "  -4.60%  "
Set-TextCell "D18" "1.977.82"
Set-TextCell "E18" "  -3.94%  "
Set-TextCell "D19" "36.178.78"
Set-TextCell "E19" "  -3.37%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D20" "0.0₃0877"
Set-TextCell "E20" "  -0.15%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D21" "70.13"
Set-TextCell "E21" "  -4.87%  "
Set-TextCell "D22" "5.27"
Set-TextCell "E22" "  -4.19%  "
Set-TextCell "D23" "233.56"
Set-TextCell "E23" "  -2.90%  "
Set-TextCell "E24" "  -0.05%  "
Set-TextCell "E25" "  -6.25%  "
Set-TextCell "E26" "  -6.85%  "
Set-TextCell "D27" "9.87"
Set-TextCell "E27" "  -1.63%  "
Set-TextCell "D28" "162.90"
Set-TextCell "E28" "  +0.24%  "
Set-TextCell "D29" "19.74"
Set-TextCell "E29" "  -2.20%  "
Set-TextCell "D30" "0.131"
Set-TextCell "E30" "  -5.45%  "
Set-TextCell "D31" "0.119"
Set-TextCell "E31" "  -3.05%  "
Set-TextCell "E32" "  -4.48%  "
Set-TextCell "D33" "4.86"
Set-TextCell "E33" "  -7.30%  "
Set-TextCell "D34" "0.0679"
Set-TextCell "E34" "  +7.03%  "
Set-TextCell "D35" "4.38"
Set-TextCell "E35" "  -7.05%  "
Set-TextCell "D36" "6.17"
Set-TextCell "E36" "  -1.72%  "
Set-TextCell "E37" "  +0.18%  "
Set-TextCell "D38" "1.81"
Set-TextCell "E38" "  -1.03%  "
Set-TextCell "D39" "2.23"
Set-TextCell "E39" "  -8.40%  "
Set-TextCell "E40" "  -5.86%  "
Set-TextCell "D41" "1.22"
Set-TextCell "E41" "  -3.99%  "
Set-TextCell "E42" "  -7.81%  "
Set-TextCell "D43" "2.89"
Set-TextCell "E43" "  -5.54%  "
Set-TextCell "D44" "0.0212"
Set-TextCell "E44" "  -3.84%  "
Set-TextCell "D45" "1.08"
Set-TextCell "E45" "  -6.10%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D46" "91.75"
Set-TextCell "E46" "  -5.71%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D47" "16.09"
Set-TextCell "E47" "  -12.17%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D48" "7.44"
Set-TextCell "E48" "  -7.06%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D49" "1.358.87"
Set-TextCell "E49" "  -4.18%  "
Set-TextCell "E50" "  -4.64%  "
Set-TextCell "D51" "44.77"
Set-TextCell "E51" "  -5.29%  "
